## FPIEBP workbook update — "updated 4.0 files and mdl"
##
## Semantic changes:
##  - About!C1: source date bumped from 2024-01-03 to 2024-03-28
##  - FPIEBP!B3:D3 (natural gas priorities): reordered from 3,2,1 to 1,3,2
##  - Cosmetic: last-used cell selection on the FPIEBP sheet, and the
##    scroll position on the About sheet.

$wb = $excel.ActiveWorkbook

$wsAbout  = $wb.Worksheets.Item("About")
$wsFpiebp = $wb.Worksheets.Item("FPIEBP")

# --- Data edits -----------------------------------------------------------

# Source date on the About tab.
$wsAbout.Range("C1").Value = 45379

# Natural gas production/imports/exports priority ordering.
$wsFpiebp.Range("B3").Value = 1
$wsFpiebp.Range("C3").Value = 3
$wsFpiebp.Range("D3").Value = 2

# --- View / selection state -------------------------------------------

# Scroll the About sheet so row 6 is at the top, and leave FPIEBP's
# selection on E3 (matching the saved cursor position in the workbook).
$wsFpiebp.Activate()
$wsFpiebp.Range("E3").Select()

$wsAbout.Activate()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1

$wsFpiebp.Activate()
$wsFpiebp.Range("E3").Select()
